$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of ICD9-CM code (column A) -> category (column C)
$map = @{
    311  = "Respiratorio"
    3129 = "Respiratorio"
    3891 = "Sangue"
    3893 = "Sangue"
    3894 = "Sangue"
    3895 = "Sangue"
    598  = "Urinario"
    5794 = "Urinario"
    8607 = "Sangue"
    8622 = "Ferita"
    8628 = "Ferita"
    8962 = "Sangue"
    8964 = "Sangue"
    9604 = "Respiratorio"
    9605 = "Respiratorio"
    9670 = "Respiratorio"
    9671 = "Respiratorio"
    9672 = "Respiratorio"
}

# Extend the Excel table (ListObject) to include the new column C first, so the
# ListColumn name picks up the header text we set afterwards.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C20"))

# Fill data rows first (matches shared-string creation order observed in target file).
# Row 3 is written before the "Colonna2" header so new unique strings land in the
# same order Excel produced them: Respiratorio, Colonna2, Sangue, Urinario, Ferita.
$code3 = $ws.Cells.Item(3, 1).Value2
$ws.Cells.Item(3, 3).Value2 = $map[[int]$code3]

# Header for new column C
$ws.Range("C1").Value2 = "Colonna2"

for ($r = 4; $r -le 20; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($map.ContainsKey([int]$code)) {
        $ws.Cells.Item($r, 3).Value2 = $map[[int]$code]
    }
}

# Widen column B slightly (matches the author's on-screen adjustment)
$ws.Columns.Item(2).ColumnWidth = 60.830729166666664

# Leave the cursor on C14, matching the author's final selection
$ws.Range("C14").Select() | Out-Null
